$d = $word.ActiveDocument

# 1) Update the "Curso (semestre ideal)" line to add the EB (7) entry.
$d.Content.Find.Execute(
    "Curso (semestre ideal): EQD (6), EQN (6)", $true, $false, $false, $false, $false,
    $true, 1, $false, "Curso (semestre ideal): EB (7), EQD (6), EQN (6)", 2)

# 2) Replace the sole requirement's course with the new weak requirement text.
$d.Content.Find.Execute(
    "LOQ4055 -  Quimica Inorgânica  (Requisito fraco)", $true, $false, $false, $false, $false,
    $true, 1, $false, "LOT2059 -  Química Orgânica Fundamental  (Requisito fraco)", 2)

# 3) Add a second weak requirement as its own run (with its own line break),
#    inside the same paragraph. We do this by splitting the paragraph in two
#    (which naturally gives us a fresh run for the new text), filling the new
#    paragraph's text via Find/Replace (so no accidental xml:space/rPr
#    artifacts are introduced), and then re-joining the two paragraphs by
#    deleting the paragraph mark between them - this leaves two separate
#    runs inside a single paragraph, matching Word's normal behaviour.
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "*LOT2059*") {
        $targetPara = $candidate
    }
}
$paraRange = $targetPara.Range
$splitPos = $paraRange.End - 1

$splitRange = $d.Range($splitPos, $splitPos)
$splitRange.InsertParagraphAfter()

$newParaStart = $splitPos + 1
$placeholder = $d.Range($newParaStart, $newParaStart)
$placeholder.InsertAfter("X")

$fillRange = $d.Range($newParaStart, $newParaStart + 1)
$fillRange.Find.Execute(
    "X", $false, $false, $false, $false, $false,
    $true, 1, $false,
    "LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito fraco)^l", 2)

$markRange = $d.Range($splitPos, $splitPos + 1)
$markRange.Delete()
